$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3097.182
$ws.Range("J17").Value = 3097.182
$ws.Range("L17").Value = 9291.545999999998
$ws.Range("N17").Value = -9627.545999999998

$ws.Range("H33").Value = 16469.85
$ws.Range("I33").Value = 19248.766
$ws.Range("J33").Value = 722.6667
$ws.Range("K33").Value = 19248.766
$ws.Range("L33").Value = 722.6667
$ws.Range("M33").Value = -19019.766
$ws.Range("N33").Value = -1180.6667

$ws.Range("H48").Value = 1448.1333
$ws.Range("I48").Value = 337.4
$ws.Range("K48").Value = 1012.2
$ws.Range("M48").Value = -720.1999999999999

$ws.Range("H56").Value = 1448.1333
$ws.Range("I56").Value = 337.4
$ws.Range("K56").Value = 1012.2
$ws.Range("M56").Value = -478.1999999999999

$ws.Range("H70").Value = 2838.4614
$ws.Range("I70").Value = 2400
$ws.Range("K70").Value = 7200
$ws.Range("M70").Value = -6930

$ws.Range("H73").Value = 2838.4614
$ws.Range("I73").Value = 2400
$ws.Range("K73").Value = 7200
$ws.Range("M73").Value = -6264

$ws.Range("H74").Value = 9237.916999999999
$ws.Range("I74").Value = 7018.4614
$ws.Range("K74").Value = 7018.4614
$ws.Range("M74").Value = -6082.4614

$ws.Range("H77").Value = 9237.916999999999
$ws.Range("I77").Value = 7018.4614
$ws.Range("K77").Value = 35092.307
$ws.Range("M77").Value = -30412.307

$ws.Range("H133").Value = 99416
$ws.Range("J133").Value = 99416
$ws.Range("L133").Value = 99416
$ws.Range("N133").Value = -109536

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11590.167
$ws.Range("I32").Value = 9323.736999999999
$ws.Range("K32").Value = 9323.736999999999
$ws.Range("M32").Value = -9036.736999999999

$ws.Range("H45").Value = 91228.30499999999
$ws.Range("I45").Value = 114458.39
$ws.Range("K45").Value = 114458.39
$ws.Range("M45").Value = -114081.39

$ws.Range("H61").Value = 6115.0264
$ws.Range("I61").Value = 6145.2163
$ws.Range("K61").Value = 6145.2163
$ws.Range("M61").Value = -5933.2163

$ws.Range("H132").Value = 1474.5094
$ws.Range("I132").Value = 1487.5
$ws.Range("K132").Value = 4462.5
$ws.Range("M132").Value = -1932.5

$ws.Range("H136").Value = 6115.0264
$ws.Range("I136").Value = 6145.2163
$ws.Range("K136").Value = 18435.6489
$ws.Range("M136").Value = -15885.6489

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1386.8846
$ws.Range("I86").Value = 1187.75
$ws.Range("K86").Value = 1187.75
$ws.Range("M86").Value = -64.75

$ws.Range("H89").Value = 1386.8846
$ws.Range("I89").Value = 1187.75
$ws.Range("K89").Value = 5938.75
$ws.Range("M89").Value = -322.75

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4675.643
$ws.Range("I16").Value = 3887.2727
$ws.Range("K16").Value = 3887.2727
$ws.Range("M16").Value = -3600.2727

$ws.Range("H33").Value = 1475
$ws.Range("I33").Value = 400
$ws.Range("J33").Value = 2012.5
$ws.Range("K33").Value = 400
$ws.Range("L33").Value = 2012.5
$ws.Range("M33").Value = -21
$ws.Range("N33").Value = -2770.5

$ws.Range("H36").Value = 444
$ws.Range("I36").Value = 444
$ws.Range("K36").Value = 444
$ws.Range("M36").Value = -56

$ws.Range("H40").Value = 444
$ws.Range("I40").Value = 444
$ws.Range("K40").Value = 444
$ws.Range("M40").Value = -284

$ws.Range("H58").Value = 9410.571
$ws.Range("I58").Value = 4852.5386
$ws.Range("K58").Value = 4852.5386
$ws.Range("M58").Value = -4649.5386

$ws.Range("H113").Value = 4675.643
$ws.Range("I113").Value = 3887.2727
$ws.Range("K113").Value = 3887.2727
$ws.Range("M113").Value = -1717.2727

$ws.Range("H132").Value = 4042.65
$ws.Range("I132").Value = 2543.5386
$ws.Range("J132").Value = 6826.7144
$ws.Range("K132").Value = 7630.6158
$ws.Range("L132").Value = 20480.1432
$ws.Range("M132").Value = -5100.6158
$ws.Range("N132").Value = -25540.1432

$ws.Range("H136").Value = 9410.571
$ws.Range("I136").Value = 4852.5386
$ws.Range("K136").Value = 14557.6158
$ws.Range("M136").Value = -12007.6158

$ws.Range("H141").Value = 166964
$ws.Range("J141").Value = 166964
$ws.Range("L141").Value = 166964
$ws.Range("N141").Value = -177324

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 262.72726
$ws.Range("I98").Value = 234
$ws.Range("K98").Value = 702
$ws.Range("M98").Value = 796

$ws.Range("H120").Value = 16171.5
$ws.Range("I120").Value = 14507.25
$ws.Range("J120").Value = 19500
$ws.Range("K120").Value = 43521.75
$ws.Range("L120").Value = 58500
$ws.Range("M120").Value = -38683.75
$ws.Range("N120").Value = -68176

$ws.Range("H122").Value = 6107.625
$ws.Range("J122").Value = 7642.8335
$ws.Range("L122").Value = 68785.5015
$ws.Range("N122").Value = -73685.5015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 49.3125
$ws.Range("I2").Value = 48.692307
$ws.Range("K2").Value = 48.692307
$ws.Range("M2").Value = 64.307693

$ws.Range("H122").Value = 2306
$ws.Range("I122").Value = 2295.7144
$ws.Range("K122").Value = 6887.1432
$ws.Range("M122").Value = -4437.1432

$ws.Range("H132").Value = 4015.7742
$ws.Range("I132").Value = 3602.2068
$ws.Range("K132").Value = 10806.6204
$ws.Range("M132").Value = -8276.6204

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 400
$ws.Range("I41").Value = 400
$ws.Range("K41").Value = 400
$ws.Range("M41").Value = 38

$ws.Range("H68").Value = 4935.75
$ws.Range("I68").Value = 2120
$ws.Range("K68").Value = 2120
$ws.Range("M68").Value = -1371

$ws.Range("H71").Value = 4935.75
$ws.Range("I71").Value = 2120
$ws.Range("K71").Value = 10600
$ws.Range("M71").Value = -6856

$ws.Range("H82").Value = 1232.5385
$ws.Range("I82").Value = 1195.2
$ws.Range("J82").Value = 1255.875
$ws.Range("K82").Value = 1195.2
$ws.Range("L82").Value = 1255.875
$ws.Range("M82").Value = -834.2
$ws.Range("N82").Value = -1977.875

$ws.Range("H85").Value = 1232.5385
$ws.Range("I85").Value = 1195.2
$ws.Range("J85").Value = 1255.875
$ws.Range("K85").Value = 1195.2
$ws.Range("L85").Value = 1255.875
$ws.Range("M85").Value = 52.79999999999995
$ws.Range("N85").Value = -3751.875

$ws.Range("H132").Value = 9211.581
$ws.Range("I132").Value = 9370.107
$ws.Range("K132").Value = 28110.321
$ws.Range("M132").Value = -25580.321

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7494
$ws.Range("J15").Value = 7494
$ws.Range("L15").Value = 7494
$ws.Range("N15").Value = -8070

$ws.Range("H62").Value = 3891.75
$ws.Range("J62").Value = 3725
$ws.Range("L62").Value = 3725
$ws.Range("N62").Value = -4973

$ws.Range("H65").Value = 3891.75
$ws.Range("J65").Value = 3725
$ws.Range("L65").Value = 18625
$ws.Range("N65").Value = -24865

$ws.Range("H132").Value = 6541.1377
$ws.Range("I132").Value = 5949.731
$ws.Range("K132").Value = 17849.193
$ws.Range("M132").Value = -15319.193

$ws.Range("H136").Value = 5479.35
$ws.Range("I136").Value = 4975.8237
$ws.Range("K136").Value = 14927.4711
$ws.Range("M136").Value = -12377.4711
